$d = $word.ActiveDocument

# --- 1) Turn the single empty paragraph after the "Inner issue" paragraph
#        into three empty paragraphs (do this first - it is the last edit
#        point in the document, so earlier paragraph indices are untouched).
$pTrailingGap = $d.Paragraphs(22)
[void]$pTrailingGap.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"/><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"/><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"/>")

# --- 2) Rewrite the "Inner issue" paragraph: drop the stray
#        <w:lastRenderedPageBreak/> and the _GoBack bookmark (both moved
#        earlier in the document, see step 4 below).
$pInnerIssue = $d.Paragraphs(21)
[void]$pInnerIssue.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t xml:space=`"preserve`">Inner issue: having to deal with someone she truly </w:t></w:r><w:r><w:t xml:space=`"preserve`">hates but is expected to - not just tolerate but – accept back into her life with open arms. </w:t></w:r></w:p>")

# --- 3) Delete the empty "ListParagraph" paragraph sitting just above
#        "Draft 2".
$pEmptyListPara = $d.Paragraphs(17)
$pEmptyListPara.Range.Delete()

# --- 4) Replace the empty paragraph right after the opening narrative
#        paragraph with the six new paragraphs (separator, two rewritten
#        narrative paragraphs with a relocated _GoBack bookmark, a new
#        centred editorial note, a rewritten paragraph with a relocated
#        lastRenderedPageBreak, and a trailing blank paragraph).
$pGap = $d.Paragraphs(8)
[void]$pGap.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:t>--------------------</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t xml:space=`"preserve`">Francesca </w:t></w:r><w:r><w:t xml:space=`"preserve`">Yates </w:t></w:r><w:r><w:t xml:space=`"preserve`">sat cross-legged the kitchen floor and stared at the opposite wall, looking at the remains of the once-fragrant yellow roses that she’d destroyed just 15 minutes prior. </w:t></w:r><w:r><w:t>It wasn’t that the flowers were awful or anything like that- she actually quite liked the roses and s</w:t></w:r><w:r><w:t>he’</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space=`"preserve`"> been excited when the delivery man had handed her the bouquet. </w:t></w:r><w:r><w:t>If she’d known who the delivery was from</w:t></w:r><w:r><w:t>, however,</w:t></w:r><w:r><w:t xml:space=`"preserve`"> she’d have saved herself the effort.</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t>Francesca poured herself a second</w:t></w:r><w:r><w:t xml:space=`"preserve`"> – or was it</w:t></w:r><w:r><w:t xml:space=`"preserve`"> a</w:t></w:r><w:r><w:t xml:space=`"preserve`"> third -</w:t></w:r><w:r><w:t xml:space=`"preserve`"> glass of red wine and checked her watch; it was 9pm. 13 hours until D-day. She took a sip of wine and sighed. She probably shouldn’t drinking it this quickly. Her best friend had bought the Merlot for her birthday and she’d managed to pour out half the bottle within two days. But Allie was the least of her concern right now. With a sigh, Francesca got up, grabbed a broom and got</w:t></w:r><w:r><w:t xml:space=`"preserve`"> to</w:t></w:r><w:r><w:t xml:space=`"preserve`"> sweeping up the remains of</w:t></w:r><w:r><w:t xml:space=`"preserve`"> the</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t xml:space=`"preserve`"> greeting card that she’d ripped to shreds straight after reading. She knelt down and picked up a piece that was curiously intact. ‘Love, Jessie’ it read. She took another sip of wine and got back to work. She would need to clean this place up as best she could if the step-monster was coming to stay with her.</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:t>[Exploring Francesca’s memory of her step-mother]</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t xml:space=`"preserve`">Francesca finished </w:t></w:r><w:r><w:t xml:space=`"preserve`">with the greeting card and took stock for a moment. At the mere mention of her step-mother’s name, she had managed to allow herself to be reduced to this: a half-drunk mess who destroys innocent flowers and drinks half a bottle of wine by herself in her pajamas on a Friday night. God, she was pathetic. A big part of her (the not-so-sober part) wanted to call the step-monster up and tell her to go fuck </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>herself instead of staying at Francesca’s place while she was in town but that same part of her also wanted to sit on the floor finishing the rest of the wine and eat ice-cream for the rest of the night. Francesca refus</w:t></w:r><w:r><w:t>ed to let herself sink that low. Besides she wasn’t 16 anymore and that kind of behavior wouldn’t cut it now that she was an adult.</w:t></w:r><w:r><w:t xml:space=`"preserve`"> She’d face the bitch and get this over and done with. It was only for a week. She could live with the woman for a week. Right?</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>")
